$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Determine the last used row dynamically so the change applies to the
# full data range (falls back to the known extent if detection fails).
$lastRow = 395
try {
    $ur = $ws.UsedRange
    $computed = $ur.Row() + $ur.Rows.Count() - 1
    if ($computed -gt 1) {
        $lastRow = $computed
    }
} catch {
    $lastRow = 395
}

# Column C ("Förändrad") holds the "last changed" date for every record.
# Bump it forward by one day (2023-09-12 -> 2023-09-13) on every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value()
    if ($v -ne $null) {
        $cell.Value = $v.AddDays(1)
    }
}
